$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1611.3334
$ws.Range("I19").Value = 1874
$ws.Range("J19").Value = 1283
$ws.Range("K19").Value = 1874
$ws.Range("L19").Value = 1283
$ws.Range("M19").Value = -1699
$ws.Range("N19").Value = -1633

$ws.Range("H40").Value = 1668.95
$ws.Range("I40").Value = 1653.303
$ws.Range("K40").Value = 1653.303
$ws.Range("M40").Value = -1478.303

$ws.Range("H88").Value = 1857.2
$ws.Range("J88").Value = 1806.6
$ws.Range("L88").Value = 1806.6
$ws.Range("N88").Value = -2618.6

$ws.Range("H91").Value = 1857.2
$ws.Range("J91").Value = 1806.6
$ws.Range("L91").Value = 1806.6
$ws.Range("N91").Value = -4614.6

$ws.Range("H111").Value = 27032.25
$ws.Range("I111").Value = 27822
$ws.Range("J111").Value = 24663
$ws.Range("K111").Value = 83466
$ws.Range("L111").Value = 73989
$ws.Range("M111").Value = -80399
$ws.Range("N111").Value = -80123

$ws.Range("H116").Value = 9725
$ws.Range("I116").Value = 12500.75
$ws.Range("K116").Value = 12500.75
$ws.Range("M116").Value = -9058.75

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = $null

$ws.Range("H135").Value = 7043
$ws.Range("I135").Value = 50
$ws.Range("K135").Value = 450
$ws.Range("M135").Value = 2085

$ws.Range("H136").Value = 147563
$ws.Range("J136").Value = 147563
$ws.Range("L136").Value = 147563
$ws.Range("N136").Value = -157763

$ws.Range("H137").Value = 1569.75
$ws.Range("I137").Value = 1538.1111
$ws.Range("J137").Value = 1664.6666
$ws.Range("K137").Value = 4614.3333
$ws.Range("L137").Value = 4993.9998
$ws.Range("M137").Value = -2064.3333
$ws.Range("N137").Value = -10093.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2753492.8
$ws.Range("I32").Value = 2596214.8
$ws.Range("K32").Value = 2596214.8
$ws.Range("M32").Value = -2595927.8

$ws.Range("H46").Value = 14748
$ws.Range("J46").Value = 4996
$ws.Range("L46").Value = 4996
$ws.Range("N46").Value = -5634

$ws.Range("H74").Value = 2672.5
$ws.Range("I74").Value = 2479.6
$ws.Range("J74").Value = 2994
$ws.Range("K74").Value = 2479.6
$ws.Range("L74").Value = 2994
$ws.Range("M74").Value = -1605.6
$ws.Range("N74").Value = -4742

$ws.Range("H77").Value = 2672.5
$ws.Range("I77").Value = 2479.6
$ws.Range("J77").Value = 2994
$ws.Range("K77").Value = 12398
$ws.Range("L77").Value = 14970
$ws.Range("M77").Value = -8030
$ws.Range("N77").Value = -23706

$ws.Range("H110").Value = 7401439.5
$ws.Range("I110").Value = 9251049
$ws.Range("K110").Value = 9251049
$ws.Range("M110").Value = -9249004

$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1048.6666
$ws.Range("I94").Value = 1048.6666
$ws.Range("K94").Value = 1048.6666
$ws.Range("M94").Value = -597.6666

$ws.Range("H105").Value = 2833.3333
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws.Range("H134").Value = 1801.7142
$ws.Range("I134").Value = 1801.7142
$ws.Range("K134").Value = 5405.142599999999
$ws.Range("M134").Value = -2870.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8666.666999999999
$ws.Range("I16").Value = 6000
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -5713
$ws.Range("N16").Value = -10574

$ws.Range("H31").Value = 2258.7693
$ws.Range("I31").Value = 2258.7693
$ws.Range("K31").Value = 2258.7693
$ws.Range("M31").Value = -1963.7693

$ws.Range("H34").Value = 2258.7693
$ws.Range("I34").Value = 2258.7693
$ws.Range("K34").Value = 2258.7693
$ws.Range("M34").Value = -2056.7693

$ws.Range("H86").Value = 5817.65
$ws.Range("I86").Value = 5808.7896
$ws.Range("K86").Value = 5808.7896
$ws.Range("M86").Value = -4685.7896

$ws.Range("H89").Value = 5817.65
$ws.Range("I89").Value = 5808.7896
$ws.Range("K89").Value = 29043.948
$ws.Range("M89").Value = -23427.948

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = $null

$ws.Range("H113").Value = 8666.666999999999
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -14340

$ws.Range("H122").Value = 2634.2307
$ws.Range("I122").Value = 1468.125
$ws.Range("K122").Value = 4404.375
$ws.Range("M122").Value = -1954.375

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null

$ws.Range("H132").Value = 5972.875
$ws.Range("I132").Value = 5548.857
$ws.Range("J132").Value = 8941
$ws.Range("K132").Value = 16646.571
$ws.Range("L132").Value = 26823
$ws.Range("M132").Value = -14116.571
$ws.Range("N132").Value = -31883

$ws.Range("H134").Value = 3050.75
$ws.Range("I134").Value = 2563.4375
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 7690.3125
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -5155.3125
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 545
$ws.Range("I55").Value = 240
$ws.Range("J55").Value = 850
$ws.Range("K55").Value = 720
$ws.Range("L55").Value = 2550
$ws.Range("M55").Value = -543
$ws.Range("N55").Value = -2904

$ws.Range("H104").Value = 6998
$ws.Range("J104").Value = 6998
$ws.Range("L104").Value = 20994
$ws.Range("N104").Value = -26236

$ws.Range("H122").Value = 150
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null

$ws.Range("H128").Value = 584296.25
$ws.Range("I128").Value = 584296.25
$ws.Range("K128").Value = 1752888.75
$ws.Range("M128").Value = -1747908.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1446.0625
$ws.Range("I107").Value = 661.1818
$ws.Range("J107").Value = 3172.8
$ws.Range("K107").Value = 661.1818
$ws.Range("L107").Value = 3172.8
$ws.Range("M107").Value = 1258.8182
$ws.Range("N107").Value = -7012.8

$ws.Range("H126").Value = 2666.6667
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8329.200000000001
$ws.Range("I7").Value = 8638.799999999999
$ws.Range("K7").Value = 8638.799999999999
$ws.Range("M7").Value = -8526.799999999999

$ws.Range("H40").Value = 3892.4285
$ws.Range("J40").Value = 5766.6665
$ws.Range("L40").Value = 5766.6665
$ws.Range("N40").Value = -6038.6665

$ws.Range("H46").Value = 1700
$ws.Range("I46").Value = 650.2
$ws.Range("K46").Value = 650.2
$ws.Range("M46").Value = -462.2

$ws.Range("H61").Value = 799
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null

$ws.Range("H113").Value = 799
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null

$ws.Range("H126").Value = 8329.200000000001
$ws.Range("I126").Value = 8638.799999999999
$ws.Range("K126").Value = 25916.4
$ws.Range("M126").Value = -23446.4

$ws.Range("H132").Value = 3574.647
$ws.Range("I132").Value = 2524.7273
$ws.Range("K132").Value = 7574.1819
$ws.Range("M132").Value = -5044.1819

$ws.Range("H136").Value = 1411.2222
$ws.Range("J136").Value = 2497.5
$ws.Range("L136").Value = 7492.5
$ws.Range("N136").Value = -12592.5

$ws.Range("H140").Value = 49900
$ws.Range("J140").Value = 49900
$ws.Range("L140").Value = 49900
$ws.Range("N140").Value = -60260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 19659.6
$ws.Range("J45").Value = 19659.6
$ws.Range("L45").Value = 19659.6
$ws.Range("N45").Value = -20641.6

$ws.Range("H82").Value = 32250
$ws.Range("I82").Value = 15000
$ws.Range("K82").Value = 15000
$ws.Range("M82").Value = -14617

$ws.Range("H85").Value = 32250
$ws.Range("I85").Value = 15000
$ws.Range("K85").Value = 15000
$ws.Range("M85").Value = -13674

$ws.Range("H132").Value = 1599.6
$ws.Range("I132").Value = 1587
$ws.Range("K132").Value = 4761
$ws.Range("M132").Value = -2231

$ws.Range("H136").Value = 1415.8182
$ws.Range("I136").Value = 1471.4
$ws.Range("J136").Value = 860
$ws.Range("K136").Value = 4414.200000000001
$ws.Range("L136").Value = 2580
$ws.Range("M136").Value = -1864.200000000001
$ws.Range("N136").Value = -7680
